$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AV1").Value = "23/04/2020"
$ws.Range("AW1").Value = "24/03/2020"
$ws.Range("AX1").Value = "24/04/2020"
$ws.Range("AY1").Value = "25/03/2020"
$ws.Range("AZ1").Value = "26/03/2020"
$ws.Range("BA1").Value = "27/03/2020"
$ws.Range("BB1").Value = "28/03/2020"
$ws.Range("BC1").Value = "29/03/2020"
$ws.Range("BD1").Value = "30/01/2020"
$ws.Range("BE1").Value = "30/03/2020"
$ws.Range("BF1").Value = "31/03/2020"
$ws.Range("AX2").Value = 108
$ws.Range("AT3").Value = 1
$ws.Range("AV3").Value = 4
$ws.Range("AX3").Value = 5
$ws.Range("AY3").Value = ""
$ws.Range("AZ3").Value = 1
$ws.Range("BA3").Value = 5
$ws.Range("BB3").Value = 3
$ws.Range("BC3").Value = ""
$ws.Range("BE3").Value = 1
$ws.Range("AT4").Value = 56
$ws.Range("AV4").Value = 80
$ws.Range("AW4").Value = 1
$ws.Range("AX4").Value = 62
$ws.Range("AZ4").Value = 1
$ws.Range("BB4").Value = 6
$ws.Range("BD4").Value = ""
$ws.Range("BE4").Value = 2
$ws.Range("BF4").Value = 21
$ws.Range("AV6").Value = 1
$ws.Range("BD6").Value = ""
$ws.Range("BF6").Value = 1
$ws.Range("AT7").Value = 17
$ws.Range("AV7").Value = 27
$ws.Range("AW7").Value = ""
$ws.Range("AX7").Value = 44
$ws.Range("AY7").Value = 1
$ws.Range("AZ7").Value = 3
$ws.Range("BA7").Value = 2
$ws.Range("BB7").Value = 2
$ws.Range("BC7").Value = 4
$ws.Range("BD7").Value = ""
$ws.Range("BF7").Value = 6
$ws.Range("AY8").Value = ""
$ws.Range("BA8").Value = 1
$ws.Range("BC8").Value = ""
$ws.Range("BD8").Value = ""
$ws.Range("BE8").Value = 5
$ws.Range("BF8").Value = 2
$ws.Range("AW9").Value = ""
$ws.Range("AX9").Value = ""
$ws.Range("AY9").Value = 2
$ws.Range("AZ9").Value = 3
$ws.Range("BB9").Value = 1
$ws.Range("BC9").Value = ""
$ws.Range("BD9").Value = ""
$ws.Range("BE9").Value = 1
$ws.Range("BF9").Value = 1
$ws.Range("AT10").Value = 92
$ws.Range("AV10").Value = 128
$ws.Range("AW10").Value = ""
$ws.Range("AX10").Value = ""
$ws.Range("AY10").Value = 5
$ws.Range("AZ10").Value = 4
$ws.Range("BA10").Value = 1
$ws.Range("BB10").Value = 9
$ws.Range("BC10").Value = 23
$ws.Range("BD10").Value = ""
$ws.Range("BE10").Value = 25
$ws.Range("BF10").Value = 23
$ws.Range("AW11").Value = ""
$ws.Range("AY11").Value = 3
$ws.Range("BA11").Value = ""
$ws.Range("BC11").Value = 2
$ws.Range("AT12").Value = 229
$ws.Range("AV12").Value = 217
$ws.Range("AW12").Value = 6
$ws.Range("AX12").Value = 191
$ws.Range("AY12").Value = 2
$ws.Range("AZ12").Value = 6
$ws.Range("BA12").Value = 3
$ws.Range("BB12").Value = 8
$ws.Range("BC12").Value = 8
$ws.Range("BD12").Value = ""
$ws.Range("BE12").Value = 7
$ws.Range("BF12").Value = 4
$ws.Range("AT13").Value = 9
$ws.Range("AV13").Value = 6
$ws.Range("AW13").Value = 2
$ws.Range("AX13").Value = 5
$ws.Range("AZ13").Value = 1
$ws.Range("BA13").Value = 1
$ws.Range("BB13").Value = 2
$ws.Range("BC13").Value = ""
$ws.Range("BD13").Value = ""
$ws.Range("BE13").Value = 1
$ws.Range("BF13").Value = 7
$ws.Range("AV14").Value = 1
$ws.Range("AT15").Value = 27
$ws.Range("AV15").Value = 27
$ws.Range("AW15").Value = 2
$ws.Range("AX15").Value = 20
$ws.Range("AY15").Value = 5
$ws.Range("AZ15").Value = 3
$ws.Range("BA15").Value = 6
$ws.Range("BB15").Value = 13
$ws.Range("BC15").Value = 5
$ws.Range("BD15").Value = ""
$ws.Range("BE15").Value = 11
$ws.Range("BF15").Value = 6
$ws.Range("AV16").Value = 7
$ws.Range("AX16").Value = 4
$ws.Range("BD16").Value = ""
$ws.Range("BF16").Value = 1
$ws.Range("AT17").Value = 9
$ws.Range("AV17").Value = 18
$ws.Range("AW17").Value = 8
$ws.Range("AX17").Value = 29
$ws.Range("AY17").Value = 10
$ws.Range("AZ17").Value = 4
$ws.Range("BA17").Value = 9
$ws.Range("BB17").Value = 12
$ws.Range("BC17").Value = 7
$ws.Range("BD17").Value = ""
$ws.Range("BE17").Value = 5
$ws.Range("BF17").Value = 13
$ws.Range("AT18").Value = 11
$ws.Range("AV18").Value = 10
$ws.Range("AW18").Value = 14
$ws.Range("AX18").Value = 3
$ws.Range("AY18").Value = 9
$ws.Range("AZ18").Value = 19
$ws.Range("BA18").Value = 39
$ws.Range("BB18").Value = 6
$ws.Range("BC18").Value = 20
$ws.Range("BD18").Value = 1
$ws.Range("BE18").Value = 32
$ws.Range("BF18").Value = 7
$ws.Range("AT20").Value = 35
$ws.Range("AV20").Value = 100
$ws.Range("AW20").Value = 1
$ws.Range("AX20").Value = 159
$ws.Range("AY20").Value = 8
$ws.Range("AZ20").Value = 5
$ws.Range("BA20").Value = 9
$ws.Range("BB20").Value = 10
$ws.Range("BC20").Value = ""
$ws.Range("BD20").Value = ""
$ws.Range("BE20").Value = 8
$ws.Range("BF20").Value = 19
$ws.Range("AT21").Value = 431
$ws.Range("AV21").Value = 778
$ws.Range("AW21").Value = 18
$ws.Range("AX21").Value = 390
$ws.Range("AY21").Value = 15
$ws.Range("AZ21").Value = 3
$ws.Range("BA21").Value = 28
$ws.Range("BB21").Value = 28
$ws.Range("BC21").Value = 22
$ws.Range("BD21").Value = ""
$ws.Range("BE21").Value = 17
$ws.Range("BF21").Value = 82
$ws.Range("AV22").Value = ""
$ws.Range("AW22").Value = 1
$ws.Range("AW24").Value = ""
$ws.Range("AY24").Value = 1
$ws.Range("AT25").Value = 4
$ws.Range("AV25").Value = 6
$ws.Range("AX25").Value = 5
$ws.Range("AZ25").Value = 1
$ws.Range("BD25").Value = ""
$ws.Range("BF25").Value = 1
$ws.Range("AT27").Value = 27
$ws.Range("AV27").Value = 5
$ws.Range("AW27").Value = 6
$ws.Range("AX27").Value = 15
$ws.Range("AY27").Value = 2
$ws.Range("AZ27").Value = 2
$ws.Range("BA27").Value = 5
$ws.Range("BC27").Value = ""
$ws.Range("BD27").Value = ""
$ws.Range("BE27").Value = 3
$ws.Range("BF27").Value = 1
$ws.Range("AT28").Value = 153
$ws.Range("AV28").Value = 76
$ws.Range("AW28").Value = ""
$ws.Range("AX28").Value = 44
$ws.Range("AY28").Value = 6
$ws.Range("AZ28").Value = 5
$ws.Range("BA28").Value = 7
$ws.Range("BB28").Value = 4
$ws.Range("BC28").Value = 5
$ws.Range("BD28").Value = ""
$ws.Range("BE28").Value = 20
$ws.Range("BF28").Value = 14
$ws.Range("AT29").Value = 33
$ws.Range("AV29").Value = 54
$ws.Range("AW29").Value = 6
$ws.Range("AX29").Value = 72
$ws.Range("AY29").Value = 8
$ws.Range("AZ29").Value = 3
$ws.Range("BA29").Value = 9
$ws.Range("BB29").Value = 4
$ws.Range("BC29").Value = 8
$ws.Range("BD29").Value = ""
$ws.Range("BE29").Value = 17
$ws.Range("BF29").Value = 57
$ws.Range("AT30").Value = 15
$ws.Range("AV30").Value = 27
$ws.Range("AW30").Value = 6
$ws.Range("AX30").Value = ""
$ws.Range("AY30").Value = 2
$ws.Range("AZ30").Value = 4
$ws.Range("BA30").Value = 14
$ws.Range("BB30").Value = 8
$ws.Range("BC30").Value = 3
$ws.Range("BD30").Value = ""
$ws.Range("BE30").Value = 7
$ws.Range("BF30").Value = 20
$ws.Range("AT32").Value = 112
$ws.Range("AV32").Value = 61
$ws.Range("AW32").Value = 1
$ws.Range("AX32").Value = 111
$ws.Range("AY32").Value = 3
$ws.Range("AZ32").Value = 4
$ws.Range("BB32").Value = 16
$ws.Range("BC32").Value = 7
$ws.Range("BD32").Value = ""
$ws.Range("BE32").Value = 24
$ws.Range("BF32").Value = 8
$ws.Range("AV33").Value = 1
$ws.Range("AW33").Value = ""
$ws.Range("AX33").Value = 1
$ws.Range("AY33").Value = 1
$ws.Range("AZ33").Value = ""
$ws.Range("BA33").Value = ""
$ws.Range("BB33").Value = 1
$ws.Range("BC33").Value = 1
$ws.Range("AV34").Value = 33
$ws.Range("AW34").Value = 2
$ws.Range("AX34").Value = 58
$ws.Range("AY34").Value = ""
$ws.Range("AZ34").Value = 1
$ws.Range("BA34").Value = 5
$ws.Range("BB34").Value = 3
$ws.Range("BC34").Value = 3
$ws.Range("BD34").Value = ""
$ws.Range("BE34").Value = 1
$ws.Range("BF34").Value = 15
